# Fix vague tags names
# - Column S ("Group Members to be:") had every row set to the generic
#   "Co-ed" value. Replace that with more specific "Men"/"Women" labels
#   for the rows where that is actually accurate; leave genuinely co-ed
#   rows alone.
# - Column Y ("Type of Small Group (check all that apply)") is missing
#   the "Students" and "Other" tag options from its comma-separated list;
#   add them back in for every row.
# - Excel re-wraps the taller text in column Y, so bump the row heights
#   to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column S: "Group Members to be:" -------------------------------
$ws.Range("S3").Value = "Men"
$ws.Range("S5").Value = "Men"
$ws.Range("S7").Value = "Women"
$ws.Range("S8").Value = "Women"
# S4 and S6 remain "Co-ed" (unchanged).

# --- Column Y: "Type of Small Group (check all that apply)" ---------
$oldTags = "Bible Study, Prayer, Freedom, Marriage, Finance, Outreach, Fitness/Health, Families, Fun/Hangout/Fellowship, College Students, Outdoor, Kids"
$newTags = "Bible Study, Prayer, Freedom, Marriage, Finance, Outreach, Fitness/Health, Families, Fun/Hangout/Fellowship, Students, Other, College Students, Outdoor, Kids"

foreach ($r in 3..8) {
    $cell = $ws.Cells.Item($r, 25)  # column Y
    if ($cell.Value2 -eq $oldTags) {
        $cell.Value = $newTags
    }
}

# --- Row heights for rows 3-8 (text grew, rows need to be taller) ---
foreach ($r in 3..8) {
    $ws.Rows.Item($r).RowHeight = 242.25
}
